# Applies the "Unveiling the Convergence of Science and Art" ->
# "Beyond the Classroom: The Value of Arts Education" rewrite described
# by the diff: title/author/email swap, body-paragraph rewordings, a
# trim of the "science & art" essay down to a shorter "arts education"
# essay, and a trailing blank paragraph added at the end of the body.

$d = $word.ActiveDocument

function Replace-Text([string]$old, [string]$new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Title / byline / email -------------------------------------------------
Replace-Text "Unveiling the Convergence of Science and Art" "Beyond the Classroom: The Value of Arts Education"
Replace-Text "Ethan James" "Jessica Albright"
Replace-Text "ethan" "jessica"
Replace-Text "james@emailcentral" "albright@schoolmail"
Replace-Text "xyz" "com"

# --- First body paragraph (intro) -------------------------------------------
Replace-Text "The intersection of science and art is a realm where boundaries blur, and creativity and knowledge converge" `
             "The arts, in their myriad forms, have long been recognized for their ability to enrich our lives and expand our horizons"

Replace-Text " It is a space where artistic expression and scientific inquiry dance in harmonious unison, each enriching the other" `
             " As educators, it is our duty to cultivate a comprehensive educational experience that values the arts alongside the traditional subjects"

Replace-Text " This essay delves into the fascinating landscape of this convergence, exploring how science inspires art and how art, in turn, informs science" `
             " This essay explores the invaluable role of arts education in shaping well-rounded, creative, and expressive individuals"

Replace-Text "In the tapestry of this convergence, art finds its muse in the wonders of the natural world and the intricacies of scientific phenomena" `
             "In a world increasingly dominated by technology and quantitative reasoning, the arts offer a refuge for the imaginative mind"

Replace-Text " Artists draw inspiration from the colors of the aurora borealis, the patterns of snowflakes, and the intricate geometry of DNA" `
             " Through engagement with visual arts, music, theater, and dance, students can cultivate their unique perspectives and find creative outlets for self-expression"

Replace-Text " These natural marvels become the raw materials of artistic creation, transformed into paintings, sculptures, and installations that capture the essence of scientific concepts" `
             " The arts provide a safe haven for exploration, where students can experiment with different forms of expression, allowing them to develop their own artistic voices"

Replace-Text "Conversely, science also finds a muse in art" "Furthermore, the arts foster essential critical thinking and problem-solving skills"

Replace-Text " The creative process inherent in art encourages scientists to think outside the conventional boundaries of their disciplines" `
             " In analyzing and interpreting works of art, students develop the ability to think critically and engage in meaningful discourse"

Replace-Text " Artists' unique perspectives and methodologies prompt scientists to approach problems from novel angles, leading to breakthroughs and innovations" `
             " They learn to appreciate different perspectives and understand the nuances of communication"

# The remainder of this paragraph (science/art "reciprocal relationship" +
# the whole "dialogue" passage) collapses into a single new sentence.
# Locate the span precisely with Find (it crosses paragraph <w:br/>s, which
# a single text Find cannot span), then replace the Range in one shot so
# the trailing, untouched "." run is left alone.
$spanStart = $d.Content
$null = $spanStart.Find.Execute(" This reciprocal relationship between science and art fosters a fertile ground for interdisciplinary exploration and discovery", `
                                 $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$spanEnd = $d.Content
$null = $spanEnd.Find.Execute(" This dialogue enriches both disciplines, fostering a deeper understanding of the world and our place within it", `
                               $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$span = $d.Range($spanStart.Start, $spanEnd.End)
$span.Text = " By actively participating in the creative process, students learn to embrace challenges, think outside the box, and find innovative solutions to problems"

# --- Summary paragraph -------------------------------------------------------
Replace-Text "The convergence of science and art is a dynamic and ever-evolving realm where creativity and knowledge intertwine" `
             "In essence, arts education provides students with the tools and skills necessary to navigate the complexities of an ever-changing world"

Replace-Text " Science inspires art, providing a wealth of ideas and wonders that fuel artistic expression" `
             " It cultivates creativity, critical thinking, and empathy, while fostering a sense of personal and cultural identity"

# Again, the tail of the summary (science/art "reciprocal relationship")
# collapses into one new sentence; replace the Range directly.
$spanStart2 = $d.Content
$null = $spanStart2.Find.Execute(" Simultaneously, art informs science, challenging conventional thinking and stimulating innovative approaches to problem-solving", `
                                  $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$spanEnd2 = $d.Content
$null = $spanEnd2.Find.Execute("these disciplines fosters interdisciplinary exploration and discovery, leading to a deeper understanding of the world and our place within it", `
                                $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$span2 = $d.Range($spanStart2.Start, $spanEnd2.End)
$span2.Text = " By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society"

# --- Trailing blank paragraph -------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
